# Update the date line at the top of the worksheet.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-08-22 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-08-23 Saturday", 2)

# Update the division problems in the table. Each "problem" row (1, 5, 9, 13, 17)
# holds five unique expressions; the rows in between are blank answer rows.
# We address every cell directly by (row, column) instead of doing a global
# text Find/Replace, because several of the new values coincide with *other*
# cells' original values (e.g. "32÷6=" is both an old value and a new value),
# which would make a sequential find-and-replace corrupt already-updated cells.

$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "32÷6="
$t.Cell(1, 2).Range.Text  = "99÷2="
$t.Cell(1, 3).Range.Text  = "88÷6="
$t.Cell(1, 4).Range.Text  = "84÷7="
$t.Cell(1, 5).Range.Text  = "35÷6="

$t.Cell(5, 1).Range.Text  = "69÷7="
$t.Cell(5, 2).Range.Text  = "77÷4="
$t.Cell(5, 3).Range.Text  = "65÷4="
$t.Cell(5, 4).Range.Text  = "33÷8="
$t.Cell(5, 5).Range.Text  = "91÷7="

$t.Cell(9, 1).Range.Text  = "32÷3="
$t.Cell(9, 2).Range.Text  = "78÷5="
$t.Cell(9, 3).Range.Text  = "40÷6="
$t.Cell(9, 4).Range.Text  = "51÷8="
$t.Cell(9, 5).Range.Text  = "27÷3="

$t.Cell(13, 1).Range.Text = "71÷4="
$t.Cell(13, 2).Range.Text = "44÷6="
$t.Cell(13, 3).Range.Text = "24÷3="
$t.Cell(13, 4).Range.Text = "75÷6="
$t.Cell(13, 5).Range.Text = "29÷5="

$t.Cell(17, 1).Range.Text = "89÷4="
$t.Cell(17, 2).Range.Text = "45÷2="
$t.Cell(17, 3).Range.Text = "25÷2="
$t.Cell(17, 4).Range.Text = "60÷3="
$t.Cell(17, 5).Range.Text = "12÷4="
